# "final files for backup" — refresh the sample rows on the CDC<->Date
# conversion helper sheet with newer dates/CDC numbers, and leave the
# selection where the author last left it (C5).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data <-> CDC")
$ws.Activate()

# Row 5: new "Data" sample date and its matching "CDC" sample number.
# (C5/F5 hold formulas and recompute automatically from these inputs.)
$ws.Range("B5").Value = 44713
$ws.Range("E5").Value = 7599

# Row 7: new "Data" sample date (C7 holds a shared formula and
# recomputes automatically).
$ws.Range("B7").Value = 44708

# Move the active selection to C5 (was B6).
$ws.Range("C5").Select()

# Best-effort: restore the minimized/positioned window chrome recorded in
# the saved file (xWindow/yWindow/windowWidth/windowHeight/minimized on
# bookViews/workbookView). Harmless no-op on hosts that don't persist raw
# window geometry back to the workbook XML.
$win = $excel.ActiveWindow
if ($win) {
    $win.Left = 5385
    $win.Top = 4335
    $win.Width = 21600
    $win.Height = 11265
    $win.WindowState = -4140
}
